# Updates the "Översikt OKÄNT" sheet:
#  - insert a brand-new case (A 27984-2024) as the new row 2, pushing every
#    existing data row down by one
#  - bump the "Förändrad" (column C) timestamp for every existing row from
#    45476 to 45477 (the run date moved forward one day)
#  - append two brand-new cases at the bottom of the table
#    (A 27999-2024, A 28053-2024)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 2; everything currently on rows 2-29 moves
#    down to rows 3-30.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).RowHeight = 15

# ---------------------------------------------------------------------------
# 2. Populate the newly inserted row 2 with the new case's data.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "A 27984-2024"          # A - Beteckning

$ws.Cells.Item(2, 2).Value = 45476                   # B - Datum
$ws.Cells.Item(2, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(2, 3).Value = 45477                   # C - Förändrad
$ws.Cells.Item(2, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(2, 4).Value = "OKÄNT"                 # D - Län
$ws.Cells.Item(2, 5).Value = "OKÄNT"                 # E - Kommun

$ws.Cells.Item(2, 7).Value = 15.4                    # G - Area (ha)
$ws.Cells.Item(2, 8).Value = 1                        # H - Fridlysta
$ws.Cells.Item(2, 9).Value = 0                        # I - Signalarter
$ws.Cells.Item(2, 10).Value = 0                       # J - NT
$ws.Cells.Item(2, 11).Value = 0                       # K - VU
$ws.Cells.Item(2, 12).Value = 0                       # L - EN
$ws.Cells.Item(2, 13).Value = 0                       # M - CR
$ws.Cells.Item(2, 14).Value = 0                       # N - RE
$ws.Cells.Item(2, 15).Value = 0                       # O - Rödlistade
$ws.Cells.Item(2, 16).Value = 0                       # P - Hotade
$ws.Cells.Item(2, 17).Value = 1                       # Q - Alla arter

$ws.Cells.Item(2, 18).Value = "Tjäder"                # R - Artnamn
$ws.Cells.Item(2, 18).WrapText = $true

$ws.Cells.Item(2, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/artfynd/A 27984-2024 artfynd.xlsx", "A 27984-2024")'          # S - Artfyndslänk
$ws.Cells.Item(2, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/kartor/A 27984-2024 karta.png", "A 27984-2024")'              # T - Kartlänk
$ws.Cells.Item(2, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomål/A 27984-2024 FSC-klagomål.docx", "A 27984-2024")'    # V - Klagomålslänk
$ws.Cells.Item(2, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/klagomålsmail/A 27984-2024 FSC-klagomål mail.docx", "A 27984-2024")' # W - Klagomålsmaillänk
$ws.Cells.Item(2, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsyn/A 27984-2024 tillsynsbegäran.docx", "A 27984-2024")'  # X - Tillsynsbegäranslänk
$ws.Cells.Item(2, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/tillsynsmail/A 27984-2024 tillsynsbegäran mail.docx", "A 27984-2024")' # Y - Tillsynsbegäransmaillänk
$ws.Cells.Item(2, 26).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SVERIGE/Logging_0000/Logging_0000/fåglar/A 27984-2024 prioriterade fågelarter.docx", "A 27984-2024")'  # Z - Fågeltillsynsbegäranslänk

# ---------------------------------------------------------------------------
# 3. Every pre-existing data row (now rows 3-30) had its "Förändrad" value
#    bumped from 45476 to 45477.
# ---------------------------------------------------------------------------
for ($r = 3; $r -le 30; $r++) {
    $ws.Cells.Item($r, 3).Value = 45477
}

# ---------------------------------------------------------------------------
# 4. Append two brand-new cases at the bottom of the table (rows 31 and 32).
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 31; A = "A 27999-2024"; B = 45476; C = 45477; G = 0.7 },
    @{ Row = 32; A = "A 28053-2024"; B = 45476; C = 45477; G = 3 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    $ws.Cells.Item($r, 1).Value = $nr.A                # A - Beteckning

    $ws.Cells.Item($r, 2).Value = $nr.B                # B - Datum
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 3).Value = $nr.C                # C - Förändrad
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = "OKÄNT"              # D - Län
    $ws.Cells.Item($r, 5).Value = "OKÄNT"              # E - Kommun

    $ws.Cells.Item($r, 7).Value = $nr.G                # G - Area (ha)
    $ws.Cells.Item($r, 8).Value = 0                    # H - Fridlysta
    $ws.Cells.Item($r, 9).Value = 0                    # I - Signalarter
    $ws.Cells.Item($r, 10).Value = 0                   # J - NT
    $ws.Cells.Item($r, 11).Value = 0                   # K - VU
    $ws.Cells.Item($r, 12).Value = 0                   # L - EN
    $ws.Cells.Item($r, 13).Value = 0                   # M - CR
    $ws.Cells.Item($r, 14).Value = 0                   # N - RE
    $ws.Cells.Item($r, 15).Value = 0                   # O - Rödlistade
    $ws.Cells.Item($r, 16).Value = 0                   # P - Hotade
    $ws.Cells.Item($r, 17).Value = 0                   # Q - Alla arter

    $ws.Cells.Item($r, 18).WrapText = $true            # R - Artnamn (empty, wrap style only)
}

Write-Output "Sheet updated: new case inserted at row 2, Forandrad dates refreshed, 2 new cases appended."
